# Scheduled-runner refresh of market-price derived columns (H:N) across
# the per-job Leve sheets (ALC/ARM/BSM/CRP/LTW/WVR). Only raw cached
# values are touched - no formulas/styles exist on this data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1666866.6
$ws.Range("I9").Value = 400
$ws.Range("K9").Value = 400
$ws.Range("M9").Value = -231

$ws.Range("H40").Value = 5097.5557
$ws.Range("J40").Value = 3399.5
$ws.Range("L40").Value = 3399.5
$ws.Range("N40").Value = -3749.5

$ws.Range("H138").Value = 5294.54
$ws.Range("J138").Value = 5649.5386
$ws.Range("L138").Value = 16948.6158
$ws.Range("N138").Value = -27228.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22897.705
$ws.Range("I32").Value = 13598.255
$ws.Range("J32").Value = 50796.06
$ws.Range("K32").Value = 13598.255
$ws.Range("L32").Value = 50796.06
$ws.Range("M32").Value = -13311.255
$ws.Range("N32").Value = -51370.06

$ws.Range("H88").Value = 4632.6665
$ws.Range("I88").Value = 3600
$ws.Range("K88").Value = 3600
$ws.Range("M88").Value = -3194

$ws.Range("H91").Value = 4632.6665
$ws.Range("I91").Value = 3600
$ws.Range("K91").Value = 3600
$ws.Range("M91").Value = -2196

$ws.Range("H97").Value = 2205.4614
$ws.Range("I97").Value = 803.75
$ws.Range("J97").Value = 2828.4443
$ws.Range("K97").Value = 803.75
$ws.Range("L97").Value = 2828.4443
$ws.Range("M97").Value = -307.75
$ws.Range("N97").Value = -3820.4443

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4905
$ws.Range("I86").Value = 4905
$ws.Range("K86").Value = 4905
$ws.Range("M86").Value = -3782

$ws.Range("H89").Value = 4905
$ws.Range("I89").Value = 4905
$ws.Range("K89").Value = 24525
$ws.Range("M89").Value = -18909

$ws.Range("H134").Value = 3575.6956
$ws.Range("I134").Value = 2219.2415
$ws.Range("K134").Value = 6657.7245
$ws.Range("M134").Value = -4122.7245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3144.8235
$ws.Range("I31").Value = 2857.0667
$ws.Range("J31").Value = 5303
$ws.Range("K31").Value = 2857.0667
$ws.Range("L31").Value = 5303
$ws.Range("M31").Value = -2562.0667
$ws.Range("N31").Value = -5893

$ws.Range("H34").Value = 3144.8235
$ws.Range("I34").Value = 2857.0667
$ws.Range("J34").Value = 5303
$ws.Range("K34").Value = 2857.0667
$ws.Range("L34").Value = 5303
$ws.Range("M34").Value = -2655.0667
$ws.Range("N34").Value = -5707

$ws.Range("H58").Value = 6883.353
$ws.Range("I58").Value = 6126.0625
$ws.Range("J58").Value = 19000
$ws.Range("K58").Value = 6126.0625
$ws.Range("L58").Value = 19000
$ws.Range("M58").Value = -5923.0625
$ws.Range("N58").Value = -19406

$ws.Range("H62").Value = 10854.923
$ws.Range("I62").Value = 11813.889
$ws.Range("K62").Value = 11813.889
$ws.Range("M62").Value = -11189.889

$ws.Range("H65").Value = 10854.923
$ws.Range("I65").Value = 11813.889
$ws.Range("K65").Value = 59069.44499999999
$ws.Range("M65").Value = -55949.44499999999

$ws.Range("H68").Value = 63614.4
$ws.Range("I68").Value = 60793.777
$ws.Range("J68").Value = 89000
$ws.Range("K68").Value = 60793.777
$ws.Range("L68").Value = 89000
$ws.Range("M68").Value = -60044.777
$ws.Range("N68").Value = -90498

$ws.Range("H71").Value = 63614.4
$ws.Range("I71").Value = 60793.777
$ws.Range("J71").Value = 89000
$ws.Range("K71").Value = 182381.331
$ws.Range("L71").Value = 267000
$ws.Range("M71").Value = -178637.331
$ws.Range("N71").Value = -274488

$ws.Range("H95").Value = 12907.286
$ws.Range("J95").Value = 12907.286
$ws.Range("L95").Value = 12907.286
$ws.Range("N95").Value = -18399.286

$ws.Range("H136").Value = 6883.353
$ws.Range("I136").Value = 6126.0625
$ws.Range("J136").Value = 19000
$ws.Range("K136").Value = 18378.1875
$ws.Range("L136").Value = 57000
$ws.Range("M136").Value = -15828.1875
$ws.Range("N136").Value = -62100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5250.9414
$ws.Range("I40").Value = 5240.643
$ws.Range("J40").Value = 5299
$ws.Range("K40").Value = 5240.643
$ws.Range("L40").Value = 5299
$ws.Range("M40").Value = -5104.643
$ws.Range("N40").Value = -5571

$ws.Range("H46").Value = 3488.1765
$ws.Range("J46").Value = 3707.6924
$ws.Range("L46").Value = 3707.6924
$ws.Range("N46").Value = -4083.6924

$ws.Range("H82").Value = 2525.1667
$ws.Range("I82").Value = 2277.6667
$ws.Range("K82").Value = 2277.6667
$ws.Range("M82").Value = -1916.6667

$ws.Range("H85").Value = 2525.1667
$ws.Range("I85").Value = 2277.6667
$ws.Range("K85").Value = 2277.6667
$ws.Range("M85").Value = -1029.6667

$ws.Range("H93").Value = 4427.8887
$ws.Range("I93").Value = 4235.2856
$ws.Range("K93").Value = 4235.2856
$ws.Range("M93").Value = -2987.2856

$ws.Range("H100").Value = 6045.4546
$ws.Range("I100").Value = 3789.5715
$ws.Range("K100").Value = 3789.5715
$ws.Range("M100").Value = -3248.5715

$ws.Range("H122").Value = 4192.4053
$ws.Range("I122").Value = 2748.5334
$ws.Range("J122").Value = 5176.864
$ws.Range("K122").Value = 8245.600199999999
$ws.Range("L122").Value = 15530.592
$ws.Range("M122").Value = -5795.600199999999
$ws.Range("N122").Value = -20430.592

$ws.Range("H132").Value = 4099.6
$ws.Range("I132").Value = 4155.4443
$ws.Range("J132").Value = 4053.9092
$ws.Range("K132").Value = 12466.3329
$ws.Range("L132").Value = 12161.7276
$ws.Range("M132").Value = -9936.332900000001
$ws.Range("N132").Value = -17221.7276

$ws.Range("H136").Value = 6219.3335
$ws.Range("I136").Value = 6663.2
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 19989.6
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -17439.6
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 25000
$ws.Range("I40").Value = 20000
$ws.Range("J40").Value = 30000
$ws.Range("K40").Value = 20000
$ws.Range("L40").Value = 30000
$ws.Range("M40").Value = -19851
$ws.Range("N40").Value = -30298
